$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcomes")

# New data for rows 4-12 (rows 2-3 remain unchanged)
$data = @(
    @(934, "Persons with heart failure", 0),
    @(938, "Hospitalization with heart failure events", 30),
    @(965, "3-point MACE", 180),
    @(967, "4-point MACE", 180),
    @(1081, "Acute Myocardial Infarction including its complications", 365),
    @(1088, "Deep Vein Thrombosis (DVT)", 365),
    @(1090, "Pulmonary Embolism", 365),
    @(1104, "RBC Transfusion (adult relevant, no auto 1yr clean window)", 0),
    @(1105, "Clostridium difficile - first episode", 0)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Remove the now-obsolete trailing rows (previously rows 13-15)
$ws.Range("A13:C15").Delete()
